# Applies the cryptos.xlsx price/volume update described in the commit
# "Updated cryptos list on Mon Jul 15 19:27:37 UTC 2024 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '63.763.01'
$ws.Range('E2').Value = '  +6.28%  '
# Row 3
$ws.Range('D3').Value = '3.415.74'
$ws.Range('E3').Value = '  +7.11%  '
# Row 4
$ws.Range('E4').Value = '  +0.00%  '
# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '579.20'
$ws.Range('E5').Value = '  +7.93%  '
# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '156.93'
$ws.Range('E6').Value = '  +8.31%  '
# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.11%  '
# Row 8
$ws.Range('D8').Value = '3.424.28'
$ws.Range('E8').Value = '  +7.23%  '
# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.535'
$ws.Range('E9').Value = '  +0.45%  '
# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '7.54'
$ws.Range('E10').Value = '  +3.01%  '
# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.122'
$ws.Range('E11').Value = '  +8.11%  '
# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.437'
$ws.Range('E12').Value = '  +1.07%  '
# Row 13
$ws.Range('D13').Value = '3.981.54'
$ws.Range('E13').Value = '  +6.49%  '
# Row 14
$ws.Range('E14').Value = '  +0.39%  '
# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.0000187'
$ws.Range('E15').Value = '  +8.58%  '
# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '27.32'
$ws.Range('E16').Value = '  +6.19%  '
# Row 17
$ws.Range('D17').Value = '63.756.25'
$ws.Range('E17').Value = '  +6.25%  '
# Row 18
$ws.Range('D18').Value = '3.408.53'
$ws.Range('E18').Value = '  +6.54%  '
# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.42'
$ws.Range('E19').Value = '  +3.12%  '
# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '14.10'
$ws.Range('E20').Value = '  +6.76%  '
# Row 21
$ws.Range('E21').Value = '  +3.57%  '
# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '389.36'
$ws.Range('E22').Value = '  +5.58%  '
# Row 23
$ws.Range('E23').Value = '  +0.19%  '
# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.538'
$ws.Range('E24').Value = '  +3.16%  '
# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '71.23'
$ws.Range('E25').Value = '  +2.48%  '
# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.62'
$ws.Range('E26').Value = '  +11.48%  '
# Row 27
$ws.Range('B27').Value = 'PEPE'
$ws.Range('C27').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.0000106'
$ws.Range('E27').Value = '  +21.86%  '
# Row 28
$ws.Range('B28').Value = 'Kaspa'
$ws.Range('C28').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.181'
$ws.Range('E28').Value = '  +6.91%  '
# Row 29
$ws.Range('E29').Value = '  +0.66%  '
# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.05'
$ws.Range('E30').Value = '  +8.38%  '
# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.57'
$ws.Range('E31').Value = '  +7.84%  '
# Row 32
$ws.Range('E32').Value = '  +13.55%  '
# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.70'
$ws.Range('E33').Value = '  +8.36%  '
# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '23.28'
$ws.Range('E34').Value = '  +3.69%  '
# Row 35
$ws.Range('E35').Value = '  -0.13%  '
# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '6.73'
$ws.Range('E36').Value = '  +2.31%  '
# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.50'
$ws.Range('E37').Value = '  +10.26%  '
# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '158.29'
$ws.Range('E38').Value = '  +0.38%  '
# Row 39
$ws.Range('B39').Value = 'EnergySwap'
$ws.Range('C39').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '27.91'
$ws.Range('E39').Value = '  +6.30%  '
# Row 40
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.88'
$ws.Range('E40').Value = '  +11.52%  '
# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0769'
$ws.Range('E41').Value = '  +8.38%  '
# Row 42
$ws.Range('D42').Value = '2.950.52'
$ws.Range('E42').Value = '  +5.99%  '
# Row 43
$ws.Range('E43').Value = '  +5.41%  '
# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.767'
$ws.Range('E44').Value = '  +6.70%  '
# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '41.47'
$ws.Range('E45').Value = '  +3.89%  '
# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '4.35'
$ws.Range('E46').Value = '  +3.37%  '
# Row 47
$ws.Range('E47').Value = '  +9.50%  '
# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '22.58'
$ws.Range('E48').Value = '  +9.91%  '
# Row 49
$ws.Range('D49').Value = '3.444.42'
$ws.Range('E49').Value = '  +6.66%  '
# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '6.36'
$ws.Range('E50').Value = '  +3.43%  '
# Row 51
$ws.Range('E51').Value = '  -1.87%  '
